# HFP_john+sally.xlsx update: roll the projection years forward by one year
# for John and Sally, and add a "year" column to the Fixed Assets table.

$wb = $excel.ActiveWorkbook

# --- John: shift years 2020-2055 -> 2021-2056 (rows 2-37) -----------------
$ws1 = $wb.Worksheets.Item("John")
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 1
}

# --- Sally: shift years 2020-2058 -> 2021-2059 (rows 2-40) -----------------
$ws2 = $wb.Worksheets.Item("Sally")
for ($r = 2; $r -le 40; $r++) {
    $cell = $ws2.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 1
}

# --- Fixed Assets: insert a new "year" column before the old column D ------
$ws4 = $wb.Worksheets.Item("Fixed Assets")
$ws4.Range("D1").EntireColumn.Insert()
$ws4.Range("D1").Value = "year"
